$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used as a staging area so numeric-looking strings can be
# pasted as literal text (matching the source inlineStr cells) without
# bumping the destination cell style index (PasteSpecial values-only
# preserves the destination cell format/style).
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"

$ws.Range("D2").Value = "42.531.91"
$ws.Range("E2").Value = "  -2.73%  "
$ws.Range("D3").Value = "2.338.42"
$ws.Range("E3").Value = "  -3.80%  "
$ws.Range("E4").Value = "  -0.19%  "
$helper.Value = "319.17"
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -1.03%  "
$helper.Value = "102.57"
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -3.18%  "
$helper.Value = "0.632"
$helper.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  -1.77%  "
$ws.Range("E8").Value = "  +0.07%  "
$helper.Value = "0.606"
$helper.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  -7.62%  "
$helper.Value = "38.97"
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  -8.35%  "
$helper.Value = "0.0914"
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -3.23%  "
$helper.Value = "8.23"
$helper.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -6.89%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$helper.Value = "0.105"
$helper.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$helper.Value = "0.970"
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -5.48%  "
$helper.Value = "15.73"
$helper.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  -9.75%  "
$ws.Range("D16").Value = "2.692.83"
$ws.Range("E16").Value = "  -3.89%  "
$ws.Range("D17").Value = "2.343.55"
$ws.Range("E17").Value = "  -3.65%  "
$ws.Range("D18").Value = "42.500.97"
$ws.Range("E18").Value = "  -2.81%  "
$helper.Value = "7.64"
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +3.67%  "
$helper.Value = "0.0000105"
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  -4.34%  "
$helper.Value = "75.54"
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -0.48%  "
$helper.Value = "3.52"
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  +1.18%  "
$helper.Value = "262.14"
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -1.09%  "
$helper.Value = "2.27"
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -8.19%  "
$helper.Value = "9.66"
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -0.30%  "
$helper.Value = "0.998"
$helper.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  -0.17%  "
$helper.Value = "11.24"
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -7.42%  "
$helper.Value = "22.92"
$helper.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  -1.25%  "
$ws.Range("E29").Value = "  -1.46%  "
$helper.Value = "173.49"
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -2.66%  "
$helper.Value = "3.01"
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -7.11%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$helper.Value = "35.00"
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -9.74%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$helper.Value = "0.0876"
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -6.60%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$helper.Value = "5.99"
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -0.14%  "
$helper.Value = "0.130"
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -1.90%  "
$helper.Value = "0.111"
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  +2.68%  "
$helper.Value = "4.48"
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -9.21%  "
$helper.Value = "0.0350"
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -6.18%  "
$helper.Value = "3.69"
$helper.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -9.24%  "
$helper.Value = "2.65"
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -8.22%  "
$helper.Value = "1.44"
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -10.87%  "
$helper.Value = "0.228"
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -2.36%  "
$helper.Value = "68.65"
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("E44").Value = "  -0.19%  "
$helper.Value = "114.06"
$helper.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -5.76%  "
$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$helper.Value = "5.47"
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -4.95%  "
$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$helper.Value = "11.46"
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -9.00%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$helper.Value = "9.02"
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -4.96%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$helper.Value = "84.18"
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  +4.22%  "
$ws.Range("D50").Value = "1.538.06"
$ws.Range("E50").Value = "  -3.17%  "
$helper.Value = "0.0988"
$helper.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -2.78%  "

$helper.Clear()
$excel.CutCopyMode = 0
